# Auto-generated edit script: updates crypto price/volume table per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "64.549.16"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = "  +1.87%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "3.241.82"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value2 = "  -1.21%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value2 = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "601.18"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value2 = "  +0.08%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "139.01"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value2 = "  +0.54%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value2 = "  -0.10%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "3.240.41"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value2 = "  -1.31%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.518"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value2 = "  +0.96%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.146"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = "  -1.13%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "5.39"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value2 = "  -0.92%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.461"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value2 = "  +0.02%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.0000245"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value2 = "  +0.97%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "35.47"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value2 = "  +3.90%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "3.777.12"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = "  -1.23%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "0.120"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value2 = "  -1.09%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "3.244.75"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value2 = "  -1.14%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "64.514.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = "  +1.77%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "6.64"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value2 = "  -2.08%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "468.97"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value2 = "  -0.83%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "14.18"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value2 = "  +2.28%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "0.710"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "7.78"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value2 = "  -0.98%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "13.46"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value2 = "  -1.67%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "84.32"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value2 = "  -0.53%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value2 = "  +0.06%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "2.74"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value2 = "  -0.05%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value2 = "  -0.01%  "

# Row 29
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value2 = "ImmutableX"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "2.14"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value2 = "  +1.57%  "

# Row 30
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value2 = "RenderToken"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "7.97"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value2 = "  -0.50%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "6.86"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value2 = "  -2.36%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "27.92"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value2 = "  -1.46%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value2 = "  -0.87%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "2.53"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value2 = "  +2.05%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value2 = "  -2.92%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "5.97"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value2 = "  +0.55%  "

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value2 = "PEPE"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.0₃0754"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value2 = "  +4.34%  "

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value2 = "OKB"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "51.86"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value2 = "  +0.04%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.0403"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value2 = "  +1.43%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "2.79"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value2 = "  +4.37%  "

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value2 = "Cosmos"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "8.21"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value2 = "  -0.36%  "

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value2 = "Bittensor"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "403.83"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value2 = "  -5.47%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.115"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value2 = "  -2.59%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "2.839.67"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value2 = "  -8.30%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.258"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value2 = "  +0.86%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "2.18"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value2 = "  +0.57%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "129.35"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value2 = "  +1.63%  "

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value2 = "InjectiveProtocol"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "26.06"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value2 = "  +0.72%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value2 = "  -0.01%  "

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value2 = "Arweave"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "35.66"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value2 = "  -1.20%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.113"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value2 = "  -0.25%  "
